$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "c-sigma"
$ws.Range("A2").Value = "style"
$ws.Range("A2").Select() | Out-Null
